# Generate Report for Handoff
# Adds two new "dependency" rows (PNG images) to each of the three sheets,
# refreshes the handoff UUID / commit hash / timestamps used by the
# existing rows, and wires up the new hyperlinks.

$wb = $excel.ActiveWorkbook

# ---- constants -------------------------------------------------------
$oldMd   = "bef8f109-06de-4819-9d2b-9e3f49d29f6c.md"
$newMd   = "034faa53-8516-4f89-bbbf-f395d49b768e.md"

$oldZhXlf = "bef8f109-06de-4819-9d2b-9e3f49d29f6c.fa4a2a1344bd61e1b34cb91c9d4bcf5cd9de2b7f.zh-cn.xlf"
$newZhXlf = "034faa53-8516-4f89-bbbf-f395d49b768e.f14aea3aed7468ff78ab598aeb98659ad5d31584.zh-cn.xlf"

$oldDeXlf = "bef8f109-06de-4819-9d2b-9e3f49d29f6c.fa4a2a1344bd61e1b34cb91c9d4bcf5cd9de2b7f.de-de.xlf"
$newDeXlf = "034faa53-8516-4f89-bbbf-f395d49b768e.f14aea3aed7468ff78ab598aeb98659ad5d31584.de-de.xlf"

$png1        = "9ef0d3c9-ef23-44ca-af40-50473fae6633.png"
$png2        = "f0140d12-a977-4aad-9aef-4db836cd13ed.png"
$png1ZhTarget = "71d4c9a30414d93a5024f4ac6a22e9474a27bbe5.png"
$png2ZhTarget = "2e8944b3fd561860cac9c7a0f41dbe8be0c89a03.png"
$png1DeTarget = "71d4c9a30414d93a5024f4ac6a22e9474a27bbe5.png"
$png2DeTarget = "2e8944b3fd561860cac9c7a0f41dbe8be0c89a03.png"

$dependencyFrom = "e2e\034faa53-8516-4f89-bbbf-f395d49b768e.md"

$overviewDate = "2016-03-25 01:22:49"
$zhDate       = "2016-03-25 01:22:45"
$deDate       = "2016-03-25 01:22:49"
$epoch        = "0001-01-01 00:00:00"

$hyperlinkColor = 15570276   # RGB(100,149,237) == FF6495ED, matches the workbook's HyperLink style
$dateFormat     = "yyyy-mm-dd HH:mm:ss"

function Style-AsHyperlink($range) {
    $range.Font.Underline = 2
    $range.Font.Color = $hyperlinkColor
}

function Add-Link($ws, $cellRef, $target, $displayText) {
    $range = $ws.Range($cellRef)
    if ($range.Hyperlinks.Count -gt 0) {
        $range.Hyperlinks.Delete()
    }
    $ws.Hyperlinks.Add($range, $target, "", "", $displayText) | Out-Null
    Style-AsHyperlink $range
}

# =======================================================================
# Sheet "Overview"
# =======================================================================
$ov = $wb.Worksheets.Item("Overview")

# refresh row 2 (existing handoff)
Add-Link $ov "A2" "https://github.com/OpenLocalizationTest/oltest/blob/034faa53-8516-4f89-bbbf-f395d49b768e/e2e/$newMd" $newMd
$ov.Range("B2").Value = "Ready for handoff"
$ov.Range("C2").Value = "Ready for handoff"
$ov.Range("D2").Value = $overviewDate
$ov.Range("D2").NumberFormat = $dateFormat

# new row 3
Add-Link $ov "A3" "https://github.com/OpenLocalizationTest/oltest/blob/034faa53-8516-4f89-bbbf-f395d49b768e/e2e/$png1" $png1
$ov.Range("B3").Value = "Ready for handoff"
$ov.Range("C3").Value = "Ready for handoff"
$ov.Range("D3").Value = $overviewDate
$ov.Range("D3").NumberFormat = $dateFormat

# new row 4
Add-Link $ov "A4" "https://github.com/OpenLocalizationTest/oltest/blob/034faa53-8516-4f89-bbbf-f395d49b768e/e2e/$png2" $png2
$ov.Range("B4").Value = "Ready for handoff"
$ov.Range("C4").Value = "Ready for handoff"
$ov.Range("D4").Value = $overviewDate
$ov.Range("D4").NumberFormat = $dateFormat

# =======================================================================
# Sheet "zh-cn"
# =======================================================================
$zh = $wb.Worksheets.Item("zh-cn")

# refresh row 2
Add-Link $zh "A2" "https://github.com/OpenLocalizationTest/oltest/blob/034faa53-8516-4f89-bbbf-f395d49b768e/e2e/$newMd" $newMd
$zh.Range("B2").Value = ".md"
$zh.Range("C2").Value = "Ready for handoff"
Add-Link $zh "D2" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/53f8a75eee818d1f96147fb2c495fd51fee67c03/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$newZhXlf" $newZhXlf
$zh.Range("E2").Value = $zhDate
$zh.Range("E2").NumberFormat = $dateFormat
$zh.Range("H2").Value = $epoch
$zh.Range("H2").NumberFormat = $dateFormat
$zh.Range("J2").Value = "Include"

# new row 3 (png dependency)
Add-Link $zh "A3" "https://github.com/OpenLocalizationTest/oltest/blob/034faa53-8516-4f89-bbbf-f395d49b768e/e2e/$png1" $png1
$zh.Range("B3").Value = ".png"
$zh.Range("C3").Value = "Ready for handoff"
Add-Link $zh "D3" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/53f8a75eee818d1f96147fb2c495fd51fee67c03/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$png1ZhTarget" $png1ZhTarget
$zh.Range("E3").Value = $zhDate
$zh.Range("E3").NumberFormat = $dateFormat
$zh.Range("H3").Value = $epoch
$zh.Range("H3").NumberFormat = $dateFormat
$zh.Range("J3").Value = "IsDependency"
$zh.Range("K3").Value = $dependencyFrom

# new row 4 (png dependency)
Add-Link $zh "A4" "https://github.com/OpenLocalizationTest/oltest/blob/034faa53-8516-4f89-bbbf-f395d49b768e/e2e/$png2" $png2
$zh.Range("B4").Value = ".png"
$zh.Range("C4").Value = "Ready for handoff"
Add-Link $zh "D4" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/53f8a75eee818d1f96147fb2c495fd51fee67c03/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$png2ZhTarget" $png2ZhTarget
$zh.Range("E4").Value = $zhDate
$zh.Range("E4").NumberFormat = $dateFormat
$zh.Range("H4").Value = $epoch
$zh.Range("H4").NumberFormat = $dateFormat
$zh.Range("J4").Value = "IsDependency"
$zh.Range("K4").Value = $dependencyFrom

# =======================================================================
# Sheet "de-de"
# =======================================================================
$de = $wb.Worksheets.Item("de-de")

# refresh row 2
Add-Link $de "A2" "https://github.com/OpenLocalizationTest/oltest/blob/034faa53-8516-4f89-bbbf-f395d49b768e/e2e/$newMd" $newMd
$de.Range("B2").Value = ".md"
$de.Range("C2").Value = "Ready for handoff"
Add-Link $de "D2" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f3d55d05de7c06ab034c0ff58e4d6dbbba8fd683/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$newDeXlf" $newDeXlf
$de.Range("E2").Value = $deDate
$de.Range("E2").NumberFormat = $dateFormat
$de.Range("H2").Value = $epoch
$de.Range("H2").NumberFormat = $dateFormat
$de.Range("J2").Value = "Include"

# new row 3 (png dependency)
Add-Link $de "A3" "https://github.com/OpenLocalizationTest/oltest/blob/034faa53-8516-4f89-bbbf-f395d49b768e/e2e/$png1" $png1
$de.Range("B3").Value = ".png"
$de.Range("C3").Value = "Ready for handoff"
Add-Link $de "D3" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f3d55d05de7c06ab034c0ff58e4d6dbbba8fd683/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$png1DeTarget" $png1DeTarget
$de.Range("E3").Value = $deDate
$de.Range("E3").NumberFormat = $dateFormat
$de.Range("H3").Value = $epoch
$de.Range("H3").NumberFormat = $dateFormat
$de.Range("J3").Value = "IsDependency"
$de.Range("K3").Value = $dependencyFrom

# new row 4 (png dependency)
Add-Link $de "A4" "https://github.com/OpenLocalizationTest/oltest/blob/034faa53-8516-4f89-bbbf-f395d49b768e/e2e/$png2" $png2
$de.Range("B4").Value = ".png"
$de.Range("C4").Value = "Ready for handoff"
Add-Link $de "D4" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f3d55d05de7c06ab034c0ff58e4d6dbbba8fd683/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$png2DeTarget" $png2DeTarget
$de.Range("E4").Value = $deDate
$de.Range("E4").NumberFormat = $dateFormat
$de.Range("H4").Value = $epoch
$de.Range("H4").NumberFormat = $dateFormat
$de.Range("J4").Value = "IsDependency"
$de.Range("K4").Value = $dependencyFrom
